$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $rng = $ws.Cells.Item($row, $col)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell 2 4 '62.798.58'
Set-TextCell 2 5 '  +5.24%  '
Set-TextCell 3 4 '3.354.15'
Set-TextCell 3 5 '  +5.23%  '
Set-TextCell 4 5 '  -0.02%  '
Set-TextCell 5 4 '571.65'
Set-TextCell 5 5 '  +6.93%  '
Set-TextCell 6 4 '152.60'
Set-TextCell 6 5 '  +5.76%  '
Set-TextCell 7 5 '  -0.13%  '
Set-TextCell 8 4 '3.354.42'
Set-TextCell 8 5 '  +5.02%  '
Set-TextCell 9 4 '0.526'
Set-TextCell 9 5 '  -0.47%  '
Set-TextCell 10 4 '7.45'
Set-TextCell 10 5 '  +1.78%  '
Set-TextCell 11 4 '0.118'
Set-TextCell 11 5 '  +5.22%  '
Set-TextCell 12 4 '0.440'
Set-TextCell 12 5 '  +3.07%  '
Set-TextCell 13 4 '3.928.69'
Set-TextCell 13 5 '  +5.06%  '
Set-TextCell 14 5 '  +0.18%  '
Set-TextCell 15 5 '  +4.89%  '
Set-TextCell 16 4 '26.90'
Set-TextCell 16 5 '  +3.86%  '
Set-TextCell 17 4 '62.756.80'
Set-TextCell 17 5 '  +5.10%  '
Set-TextCell 18 4 '3.333.18'
Set-TextCell 18 5 '  +4.93%  '
Set-TextCell 19 4 '6.34'
Set-TextCell 19 5 '  +1.94%  '
Set-TextCell 20 4 '13.84'
Set-TextCell 20 5 '  +5.68%  '
Set-TextCell 21 4 '8.41'
Set-TextCell 21 5 '  +2.66%  '
Set-TextCell 22 4 '384.10'
Set-TextCell 22 5 '  +4.76%  '
Set-TextCell 23 5 '  +0.08%  '
Set-TextCell 24 4 '0.535'
Set-TextCell 24 5 '  +3.02%  '
Set-TextCell 25 4 '70.08'
Set-TextCell 25 5 '  +0.78%  '
Set-TextCell 26 4 '9.35'
Set-TextCell 26 5 '  +6.72%  '
Set-TextCell 27 5 '  +6.85%  '
Set-TextCell 28 4 '0.0₃0966'
Set-TextCell 28 5 '  +9.82%  '
Set-TextCell 29 5 '  -0.19%  '
Set-TextCell 30 5 '  +6.43%  '
Set-TextCell 31 4 '23.03'
Set-TextCell 31 5 '  +3.67%  '
Set-TextCell 32 5 '  +6.15%  '
Set-TextCell 33 4 '6.35'
Set-TextCell 33 5 '  +4.70%  '
Set-TextCell 34 5 '  +10.32%  '
Set-TextCell 35 4 '6.72'
Set-TextCell 35 5 '  +2.60%  '
Set-TextCell 36 5 '  +9.81%  '
Set-TextCell 37 4 '157.76'
Set-TextCell 37 5 '  +0.60%  '
Set-TextCell 38 4 '1.88'
Set-TextCell 38 5 '  +13.12%  '
Set-TextCell 39 4 '27.09'
Set-TextCell 39 5 '  +5.41%  '
Set-TextCell 40 5 '  +12.28%  '
Set-TextCell 41 4 '0.0736'
Set-TextCell 41 5 '  +5.83%  '
Set-TextCell 42 4 '2.778.10'
Set-TextCell 42 5 '  -0.26%  '
Set-TextCell 43 4 '40.92'
Set-TextCell 43 5 '  +4.29%  '
Set-TextCell 44 5 '  +1.89%  '
Set-TextCell 45 4 '0.745'
Set-TextCell 45 5 '  +4.84%  '
Set-TextCell 46 4 '1.04'
Set-TextCell 46 5 '  +6.35%  '
Set-TextCell 47 4 '3.395.23'
Set-TextCell 48 4 '21.94'
Set-TextCell 48 5 '  +7.87%  '
Set-TextCell 49 4 '6.34'
Set-TextCell 49 5 '  +3.70%  '
Set-TextCell 50 5 '  -0.34%  '
Set-TextCell 51 4 '291.43'
Set-TextCell 51 5 '  +10.58%  '
